# Insert two new weekly-report rows at the top of the "Tomate" block
# (rows 1195-1196), shifting all existing rows at/after 1195 down by 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1195:1196").Insert()

# ---- Row 1195 ----
$ws.Cells.Item(1195, 1).Value = 5
$ws.Cells.Item(1195, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1195, 3).Value = "Maule"
$ws.Cells.Item(1195, 4).Value = 45265
$ws.Cells.Item(1195, 5).Value = 7
$ws.Cells.Item(1195, 6).Value = 100112020
$ws.Cells.Item(1195, 7).Value = "Tomate"
$ws.Cells.Item(1195, 8).Value = "Larga vida"
$ws.Cells.Item(1195, 9).Value = "Primera"
$ws.Cells.Item(1195, 10).Value = 2500
$ws.Cells.Item(1195, 11).Value = 9000
$ws.Cells.Item(1195, 12).Value = 9000
$ws.Cells.Item(1195, 13).Value = 9000
$ws.Cells.Item(1195, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1195, 15).Value = "Región del Maule"
$ws.Cells.Item(1195, 16).Value = 500
$ws.Cells.Item(1195, 17).Value = 18
$ws.Cells.Item(1195, 18).Value = "Hortaliza"

# ---- Row 1196 ----
$ws.Cells.Item(1196, 1).Value = 5
$ws.Cells.Item(1196, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1196, 3).Value = "Maule"
$ws.Cells.Item(1196, 4).Value = 45265
$ws.Cells.Item(1196, 5).Value = 7
$ws.Cells.Item(1196, 6).Value = 100112020
$ws.Cells.Item(1196, 7).Value = "Tomate"
$ws.Cells.Item(1196, 8).Value = "Larga vida"
$ws.Cells.Item(1196, 9).Value = "Primera"
$ws.Cells.Item(1196, 10).Value = 2500
$ws.Cells.Item(1196, 11).Value = 5000
$ws.Cells.Item(1196, 12).Value = 5000
$ws.Cells.Item(1196, 13).Value = 5000
$ws.Cells.Item(1196, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(1196, 15).Value = "Región del Maule"
$ws.Cells.Item(1196, 16).Value = 500
$ws.Cells.Item(1196, 17).Value = 10
$ws.Cells.Item(1196, 18).Value = "Hortaliza"
